$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 100; this shifts existing rows 100-107
# down to 102-109, preserving all their values, formats and formulas.
$ws.Range("A100:A101").EntireRow.Insert()

# --- New row 100 ---
$ws.Cells.Item(100, 1).Value = 4
$ws.Cells.Item(100, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(100, 3).Value = "Los Lagos"
$ws.Cells.Item(100, 4).Value = 44516
$ws.Cells.Item(100, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(100, 5).Value = 10
$ws.Cells.Item(100, 6).Value = "Fruta"
$ws.Cells.Item(100, 7).Value = 100108
$ws.Cells.Item(100, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(100, 9).Value = 100108002
$ws.Cells.Item(100, 10).Value = "Mango"
$ws.Cells.Item(100, 11).Value = "Sin especificar"
$ws.Cells.Item(100, 12).Value = "Primera"
$ws.Cells.Item(100, 13).Value = 200
$ws.Cells.Item(100, 14).Value = 7500
$ws.Cells.Item(100, 15).Value = 8000
$ws.Cells.Item(100, 16).Value = 7750
$ws.Cells.Item(100, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(100, 18).Value = "Perú"
$ws.Cells.Item(100, 19).Value = 1938
$ws.Cells.Item(100, 20).Value = 4

# --- New row 101 ---
$ws.Cells.Item(101, 1).Value = 4
$ws.Cells.Item(101, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(101, 3).Value = "Los Lagos"
$ws.Cells.Item(101, 4).Value = 44516
$ws.Cells.Item(101, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(101, 5).Value = 10
$ws.Cells.Item(101, 6).Value = "Fruta"
$ws.Cells.Item(101, 7).Value = 100108
$ws.Cells.Item(101, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(101, 9).Value = 100108002
$ws.Cells.Item(101, 10).Value = "Mango"
$ws.Cells.Item(101, 11).Value = "Sin especificar"
$ws.Cells.Item(101, 12).Value = "Segunda"
$ws.Cells.Item(101, 13).Value = 100
$ws.Cells.Item(101, 14).Value = 5500
$ws.Cells.Item(101, 15).Value = 5500
$ws.Cells.Item(101, 16).Value = 5500
$ws.Cells.Item(101, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(101, 18).Value = "Perú"
$ws.Cells.Item(101, 19).Value = 1375
$ws.Cells.Item(101, 20).Value = 4
